# Adds a new "2022-Q3" sheet (cloned from the existing "2022-Q2" sheet so it
# keeps identical layout/styling) populated with the new quarter's fund
# holdings data, and updates the "总计" (summary) sheet with the new
# quarter's totals, shifting the previous rows down.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)   # "总计"
$wsQ2    = $wb.Worksheets.Item(2)   # currently "2022-Q2"

# --- 1. Create the new "2022-Q3" sheet right before the current "2022-Q2" ---
# Copying keeps all formatting/styles (header row + index column) identical.
$wsQ2.Copy($wsQ2, [System.Reflection.Missing]::Value)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# helper-free: set a text cell while preventing Excel from auto-converting
# numeric-looking strings (e.g. "009467", "1.69") into numbers, and without
# leaving a stray "quote prefix" style behind.
$wsQ3.Cells.Item(2,2).Formula = "'009467"
$wsQ3.Cells.Item(2,2).Style = "Normal"
$wsQ3.Cells.Item(2,3).Formula = "'红土创新科技创新3个月定开混合A"
$wsQ3.Cells.Item(2,3).Style = "Normal"
$wsQ3.Cells.Item(2,4).Formula = "'1.69"
$wsQ3.Cells.Item(2,4).Style = "Normal"
$wsQ3.Cells.Item(2,5).Formula = "'94.87"
$wsQ3.Cells.Item(2,5).Style = "Normal"
$wsQ3.Cells.Item(2,6).Formula = "'5.03"
$wsQ3.Cells.Item(2,6).Style = "Normal"
$wsQ3.Cells.Item(2,7).Formula = "'0.0850"
$wsQ3.Cells.Item(2,7).Style = "Normal"
$wsQ3.Cells.Item(2,8).Value = 5

$wsQ3.Cells.Item(3,2).Formula = "'013173"
$wsQ3.Cells.Item(3,2).Style = "Normal"
$wsQ3.Cells.Item(3,3).Formula = "'红土创新科技创新3个月定开混合C"
$wsQ3.Cells.Item(3,3).Style = "Normal"
$wsQ3.Cells.Item(3,4).Formula = "'0.47"
$wsQ3.Cells.Item(3,4).Style = "Normal"
$wsQ3.Cells.Item(3,5).Formula = "'94.87"
$wsQ3.Cells.Item(3,5).Style = "Normal"
$wsQ3.Cells.Item(3,6).Formula = "'5.03"
$wsQ3.Cells.Item(3,6).Style = "Normal"
$wsQ3.Cells.Item(3,7).Formula = "'0.0236"
$wsQ3.Cells.Item(3,7).Style = "Normal"
$wsQ3.Cells.Item(3,8).Value = 5

# --- 2. Update the "总计" sheet with the new quarter on top, shifting the
#        previously-existing rows down by one ---
$wsTotal.Cells.Item(2,2).Formula = "'2022-Q3"
$wsTotal.Cells.Item(2,2).Style = "Normal"
$wsTotal.Cells.Item(2,4).Value = 0.11

$wsTotal.Cells.Item(3,2).Formula = "'2022-Q2"
$wsTotal.Cells.Item(3,2).Style = "Normal"
$wsTotal.Cells.Item(3,4).Value = 0.09

$wsTotal.Cells.Item(3,1).Copy()
$wsTotal.Cells.Item(4,1).PasteSpecial(-4122)
$wsTotal.Cells.Item(4,1).Value = 2
$wsTotal.Cells.Item(4,2).Formula = "'2022-Q1"
$wsTotal.Cells.Item(4,2).Style = "Normal"
$wsTotal.Cells.Item(4,3).Value = 2
$wsTotal.Cells.Item(4,4).Value = 0.02

# --- 3. Keep the originally-selected tab ("2022-Q1") selected; it is now
#        the 4th sheet instead of the 3rd ---
$wb.Worksheets.Item(4).Activate()
